# EDV.xlsx update
# -----------------
# The "Summe" column (C) compounds the previous row by a factor of
# 2^exponent every period. The exponent is tuned down from 1.5 to 1.15
# (a much slower growth curve), which in turn changes every dependent
# value in column D ("Summe inkl. Mehrwertsteuer", which just applies
# VAT on top of C) as well as the line chart built from B:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the compounding formula in C5:C16 with the new exponent.
# (D4:D16 already reference column C, so they recalc automatically.)
for ($r = 5; $r -le 16; $r++) {
    $prevRow = $r - 1
    $ws.Cells.Item($r, 3).Formula = "=C$prevRow*2^(1.15)"
}

# The chart's log-scaled value axis had its ceiling fixed at 10,000,000
# to comfortably fit the old (much larger) curve; rescale it down to
# 550,000 to match the new, smaller maximum (~526k).
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.MaximumScale = 550000

# Move the cell selection off of C5 (where the edit was made).
$ws.Range("Q13").Select()
